$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FEB-22")
$ws.Activate()

# --- Row 3: remove the "No" value (task continues from row above) ---
$ws.Cells.Item(3, 1).ClearContents()

# --- Row 4: renumber 3 -> 2 ---
$ws.Cells.Item(4, 1).Value2 = 2

# --- Row 5: new task entry (No=3, Date=2022-02-03) ---
$ws.Cells.Item(5, 1).Value2 = 3
# B5 needs the same date number format as B4 (reuse the existing style instead of minting a new one)
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(5, 2).Value2 = 44595

# --- Row 7: new task entry (No=4, Date=2022-02-03, Application=RPA GSS) ---
$ws.Rows.Item(7).RowHeight = 28.8
$ws.Cells.Item(7, 1).Value2 = 4
$ws.Cells.Item(7, 2).Value2 = 44595
$ws.Cells.Item(7, 3).Value2 = "RPA GSS"
$ws.Cells.Item(7, 4).Value2 = "1. 5000 thousands records issue was fixed at DRS Monthly task when record count is equal to 5000 and the testing is going intermittently as the new bot is occupied by Vijay san and Nirmal san,"
$ws.Cells.Item(7, 5).Value2 = 0.9
$ws.Cells.Item(7, 6).Value2 = "WIP"

# --- Row 8: continuation comment row ---
$ws.Cells.Item(8, 4).Value2 = "2.  Implementation of Public holidays are work in progress"
$ws.Cells.Item(8, 5).Value2 = 0.4
$ws.Cells.Item(8, 6).Value2 = "WIP"

# --- Selection moved to F8 ---
$ws.Range("F8").Select()
